# Apply updated crypto price/volume figures scraped on 2024-01-24.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $rng = $ws.Range($cellRef)
    # Leading apostrophe forces Excel to treat the value as literal text
    # (prevents '40.102.76' / '22.90' style numbers from being coerced to
    # numeric values and losing their original formatting).
    $rng.Value = "'" + $text
    # Restore the default cell style so the quote-prefix flag added above
    # does not linger as a visible style change on the cell.
    $rng.Style = "Normal"
}

Set-TextValue "D2" "40.102.76"
Set-TextValue "E2" "  +1.70%  "
Set-TextValue "D3" "2.235.74"
Set-TextValue "E3" "  +1.16%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "292.59"
Set-TextValue "E5" "  -1.52%  "
Set-TextValue "D6" "87.14"
Set-TextValue "E6" "  +5.67%  "
Set-TextValue "D7" "0.516"
Set-TextValue "E7" "  +1.17%  "
Set-TextValue "E8" "  -0.07%  "
Set-TextValue "D9" "0.475"
Set-TextValue "E9" "  +1.66%  "
Set-TextValue "D10" "31.25"
Set-TextValue "E10" "  +7.54%  "
Set-TextValue "D11" "0.0789"
Set-TextValue "E11" "  +1.91%  "
Set-TextValue "D12" "47.27"
Set-TextValue "E12" "  -0.83%  "
Set-TextValue "E13" "  +1.67%  "
Set-TextValue "D14" "6.38"
Set-TextValue "E14" "  +1.58%  "
Set-TextValue "D15" "2.583.97"
Set-TextValue "E15" "  +0.92%  "
Set-TextValue "D16" "14.16"
Set-TextValue "E16" "  +0.63%  "
Set-TextValue "D17" "2.225.68"
Set-TextValue "E17" "  +0.70%  "
Set-TextValue "D18" "0.733"
Set-TextValue "E18" "  +2.68%  "
Set-TextValue "D19" "40.060.23"
Set-TextValue "E19" "  +1.92%  "
Set-TextValue "D20" "0.0₃0888"
Set-TextValue "E20" "  +1.44%  "
Set-TextValue "D21" "11.32"
Set-TextValue "E21" "  +9.24%  "
Set-TextValue "E22" "  +2.26%  "
Set-TextValue "D23" "65.87"
Set-TextValue "E23" "  +1.23%  "
Set-TextValue "D24" "236.71"
Set-TextValue "E24" "  +4.28%  "
Set-TextValue "E25" "  +0.11%  "
Set-TextValue "E26" "  +3.22%  "
Set-TextValue "E27" "  +3.43%  "
Set-TextValue "D28" "22.90"
Set-TextValue "E28" "  +1.56%  "
Set-TextValue "E29" "  +2.21%  "
Set-TextValue "D30" "9.32"
Set-TextValue "E30" "  +2.62%  "
Set-TextValue "D31" "33.26"
Set-TextValue "E31" "  +4.86%  "
Set-TextValue "D32" "151.99"
Set-TextValue "E32" "  +1.90%  "
Set-TextValue "E33" "  -0.10%  "
Set-TextValue "E34" "  +3.16%  "
Set-TextValue "D35" "0.0724"
Set-TextValue "E35" "  +4.50%  "
Set-TextValue "E36" "  +2.13%  "
Set-TextValue "D37" "16.25"
Set-TextValue "E37" "  +6.63%  "
Set-TextValue "E38" "  +8.13%  "
Set-TextValue "E39" "  +2.12%  "
Set-TextValue "E40" "  +3.62%  "
Set-TextValue "D41" "1.71"
Set-TextValue "D42" "3.82"
Set-TextValue "E42" "  +5.59%  "
Set-TextValue "D43" "2.072.73"
Set-TextValue "E43" "  +8.75%  "
Set-TextValue "D44" "18.38"
Set-TextValue "E44" "  +15.46%  "
Set-TextValue "D45" "2.13"
Set-TextValue "E45" "  +4.19%  "
Set-TextValue "D46" "0.0269"
Set-TextValue "E46" "  +4.40%  "
Set-TextValue "D47" "9.90"
Set-TextValue "E47" "  +9.93%  "
Set-TextValue "D48" "2.62"
Set-TextValue "E48" "  +0.10%  "
Set-TextValue "D49" "2.454.13"
Set-TextValue "E49" "  +1.02%  "
Set-TextValue "D50" "72.36"
Set-TextValue "E50" "  +2.05%  "
Set-TextValue "D51" "89.62"
Set-TextValue "E51" "  +2.93%  "
